$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Roraima -> Acre
$ws.Range("A2").Value = "Acre"
$ws.Range("B2").Value = "Diferença 2024/10 - 2023/10"
$ws.Range("C2").Value = 2.97

# Row 3: Mato Grosso -> Rio Grande do Norte
$ws.Range("A3").Value = "Rio Grande do Norte"
$ws.Range("B3").Value = "Diferença 2024/10 - 2023/10"
$ws.Range("C3").Value = 2.21

# Row 4: Paraíba -> Rondônia
$ws.Range("A4").Value = "Rondônia"
$ws.Range("B4").Value = "Diferença 2024/10 - 2023/10"
$ws.Range("C4").Value = 1.89

# Row 5: Rio Grande do Norte -> Tocantins
$ws.Range("A5").Value = "Tocantins"
$ws.Range("B5").Value = "Diferença 2024/10 - 2023/10"
$ws.Range("C5").Value = 1.71

# Row 6: Acre -> Roraima
$ws.Range("A6").Value = "Roraima"
$ws.Range("B6").Value = "Diferença 2024/10 - 2023/10"
$ws.Range("C6").Value = 1.61

# Row 7: Tocantins -> Pernambuco
$ws.Range("A7").Value = "Pernambuco"
$ws.Range("B7").Value = "Diferença 2024/10 - 2023/10"
$ws.Range("C7").Value = 1.58

# Row 8: Sergipe (unchanged name)
$ws.Range("B8").Value = "Diferença 2024/10 - 2023/10"
$ws.Range("C8").Value = 1.45
$ws.Range("D8").Value = "9º"

# Row 9: Brasil (unchanged name)
$ws.Range("B9").Value = "Diferença 2024/10 - 2023/10"
$ws.Range("C9").Value = 0.4

# Row 10: Nordeste (unchanged name)
$ws.Range("B10").Value = "Diferença 2024/10 - 2023/10"
$ws.Range("C10").Value = 0.91
